# Update column F (dSF) values for rows 4-23 on the active worksheet,
# per the repull/recalculation described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -2
    5  = -1
    6  = 3
    7  = 1
    8  = -1
    9  = 2
    11 = -4
    12 = 1
    13 = 2
    14 = 8
    15 = -4
    16 = 7
    18 = 2
    19 = -2
    20 = -5
    21 = 6
    22 = 4
    23 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
